$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2006.125
$ws.Range("I33").Value = 1609.5454
$ws.Range("J33").Value = 2878.6
$ws.Range("K33").Value = 1609.5454
$ws.Range("L33").Value = 2878.6
$ws.Range("M33").Value = -1380.5454
$ws.Range("N33").Value = -3336.6

$ws.Range("H53").Value = 633.28
$ws.Range("I53").Value = 828.1429000000001
$ws.Range("K53").Value = 828.1429000000001
$ws.Range("M53").Value = -191.1429000000001

$ws.Range("H96").Value = 2204.8572
$ws.Range("I96").Value = 2794.8
$ws.Range("J96").Value = 730
$ws.Range("K96").Value = 8384.400000000001
$ws.Range("L96").Value = 2190
$ws.Range("M96").Value = -7011.400000000001
$ws.Range("N96").Value = -4936

$ws.Range("H116").Value = 3411.5715
$ws.Range("I116").Value = 2816.8
$ws.Range("J116").Value = 4898.5
$ws.Range("K116").Value = 2816.8
$ws.Range("L116").Value = 4898.5
$ws.Range("M116").Value = 625.1999999999998
$ws.Range("N116").Value = -11782.5

$ws.Range("H138").Value = 3473.6765
$ws.Range("I138").Value = 3702.0667
$ws.Range("J138").Value = 3293.3684
$ws.Range("K138").Value = 11106.2001
$ws.Range("L138").Value = 9880.1052
$ws.Range("M138").Value = -5966.2001
$ws.Range("N138").Value = -20160.1052

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2499.3333
$ws.Range("J61").Value = 2499.5
$ws.Range("L61").Value = 2499.5
$ws.Range("N61").Value = -2923.5

$ws.Range("H74").Value = 1875.3684
$ws.Range("I74").Value = 1743.7059
$ws.Range("J74").Value = 2994.5
$ws.Range("K74").Value = 1743.7059
$ws.Range("L74").Value = 2994.5
$ws.Range("M74").Value = -869.7058999999999
$ws.Range("N74").Value = -4742.5

$ws.Range("H77").Value = 1875.3684
$ws.Range("I77").Value = 1743.7059
$ws.Range("J77").Value = 2994.5
$ws.Range("K77").Value = 8718.529500000001
$ws.Range("L77").Value = 14972.5
$ws.Range("M77").Value = -4350.529500000001
$ws.Range("N77").Value = -23708.5

$ws.Range("H97").Value = 1255.92
$ws.Range("I97").Value = 1226.7273
$ws.Range("J97").Value = 1470
$ws.Range("K97").Value = 1226.7273
$ws.Range("L97").Value = 1470
$ws.Range("M97").Value = -730.7273
$ws.Range("N97").Value = -2462

$ws.Range("H102").Value = 2928
$ws.Range("I102").Value = 2680.6924
$ws.Range("J102").Value = 3999.6667
$ws.Range("K102").Value = 2680.6924
$ws.Range("L102").Value = 3999.6667
$ws.Range("M102").Value = -1058.6924
$ws.Range("N102").Value = -7243.6667

$ws.Range("H132").Value = 3905.257
$ws.Range("I132").Value = 3176.1333
$ws.Range("K132").Value = 9528.3999
$ws.Range("M132").Value = -6998.3999

$ws.Range("H136").Value = 2499.3333
$ws.Range("J136").Value = 2499.5
$ws.Range("L136").Value = 7498.5
$ws.Range("N136").Value = -12598.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 13158.925
$ws.Range("I86").Value = 12578.64
$ws.Range("J86").Value = 14126.066
$ws.Range("K86").Value = 12578.64
$ws.Range("L86").Value = 14126.066
$ws.Range("M86").Value = -11455.64
$ws.Range("N86").Value = -16372.066

$ws.Range("H89").Value = 13158.925
$ws.Range("I89").Value = 12578.64
$ws.Range("J89").Value = 14126.066
$ws.Range("K89").Value = 62893.2
$ws.Range("L89").Value = 70630.33
$ws.Range("M89").Value = -57277.2
$ws.Range("N89").Value = -81862.33

$ws.Range("H94").Value = 1132.4546
$ws.Range("I94").Value = 1367
$ws.Range("J94").Value = 998.4286
$ws.Range("K94").Value = 1367
$ws.Range("L94").Value = 998.4286
$ws.Range("M94").Value = -916
$ws.Range("N94").Value = -1900.4286

$ws.Range("H105").Value = 1809.6666
$ws.Range("I105").Value = 1971.5
$ws.Range("K105").Value = 1971.5
$ws.Range("M105").Value = -224.5

$ws.Range("H106").Value = 25966.666
$ws.Range("J106").Value = 25966.666
$ws.Range("L106").Value = 25966.666
$ws.Range("N106").Value = -28490.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 15222.923
$ws.Range("J86").Value = 19973
$ws.Range("L86").Value = 19973
$ws.Range("N86").Value = -22219

$ws.Range("H89").Value = 15222.923
$ws.Range("J89").Value = 19973
$ws.Range("L89").Value = 99865
$ws.Range("N89").Value = -111097

$ws.Range("H134").Value = 1844.7693
$ws.Range("I134").Value = 1925.7778
$ws.Range("J134").Value = 1662.5
$ws.Range("K134").Value = 5777.3334
$ws.Range("L134").Value = 4987.5
$ws.Range("M134").Value = -3242.3334
$ws.Range("N134").Value = -10057.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 800008.9399999999
$ws.Range("I2").Value = 57
$ws.Range("K2").Value = 342
$ws.Range("M2").Value = -229

$ws.Range("H34").Value = 39317.43
$ws.Range("J34").Value = 42301.5
$ws.Range("L34").Value = 126904.5
$ws.Range("N34").Value = -127072.5

$ws.Range("H52").Value = 2312.75
$ws.Range("J52").Value = 2312.75
$ws.Range("L52").Value = 6938.25
$ws.Range("N52").Value = -7470.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6700
$ws.Range("I70").Value = 3400
$ws.Range("K70").Value = 3400
$ws.Range("M70").Value = -3130

$ws.Range("H73").Value = 6700
$ws.Range("I73").Value = 3400
$ws.Range("K73").Value = 3400
$ws.Range("M73").Value = -2464

$ws.Range("H132").Value = 11114581
$ws.Range("I132").Value = 3456.9048
$ws.Range("J132").Value = 37040540
$ws.Range("K132").Value = 10370.7144
$ws.Range("L132").Value = 111121620
$ws.Range("M132").Value = -7840.714399999999
$ws.Range("N132").Value = -111126680

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1557.7858
$ws.Range("I55").Value = 1285.8
$ws.Range("K55").Value = 1285.8
$ws.Range("M55").Value = -1112.8

$ws.Range("H70").Value = 30000
$ws.Range("J70").Value = 30000
$ws.Range("L70").Value = 30000
$ws.Range("N70").Value = -30540

$ws.Range("H73").Value = 30000
$ws.Range("J73").Value = 30000
$ws.Range("L73").Value = 30000
$ws.Range("N73").Value = -31872

$ws.Range("H74").Value = 25000
$ws.Range("I74").Value = 25000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 25000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -24002
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 25000
$ws.Range("I77").Value = 25000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 75000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -70008
$ws.Range("N77").ClearContents()

$ws.Range("H93").Value = 2165.6667
$ws.Range("I93").Value = 1498.5
$ws.Range("J93").Value = 3500
$ws.Range("K93").Value = 1498.5
$ws.Range("L93").Value = 3500
$ws.Range("M93").Value = -250.5
$ws.Range("N93").Value = -5996

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws.Range("H132").Value = 3043
$ws.Range("I132").Value = 2641.7144
$ws.Range("J132").Value = 3218.5625
$ws.Range("K132").Value = 7925.1432
$ws.Range("L132").Value = 9655.6875
$ws.Range("M132").Value = -5395.1432
$ws.Range("N132").Value = -14715.6875

$ws.Range("H134").Value = 125000
$ws.Range("J134").Value = 125000
$ws.Range("L134").Value = 125000
$ws.Range("N134").Value = -135140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 8339832.5
$ws.Range("J3").Value = 9749
$ws.Range("L3").Value = 9749
$ws.Range("N3").Value = -9977

$ws.Range("H11").Value = 3500.5
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 3500.5
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 3500.5
$ws.Range("N11").Value = -3784.5
$ws.Range("M11").ClearContents()

$ws.Range("H44").Value = 50000
$ws.Range("J44").Value = 50000
$ws.Range("L44").Value = 50000
$ws.Range("N44").Value = -51108

$ws.Range("H117").Value = 77500
$ws.Range("J117").Value = 77500
$ws.Range("L117").Value = 77500
$ws.Range("N117").Value = -86678

$ws.Range("H132").Value = 43488456
$ws.Range("I132").Value = 11444.7
$ws.Range("K132").Value = 34334.10000000001
$ws.Range("M132").Value = -31804.10000000001
